$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet already contains a "Unterliste" (sub-list) hierarchy ending at
# row 26 (A26/B26 = "dritte Unterliste" / "Spezialzeichen ..."). The edit
# appends three more rows (27-29) that nest additional list levels using
# the same A/B values, plus new D/E/F values "sehr" / "tief" /
# "verschachtelt!" building up one column at a time.

# Row 27
$ws.Range("A27").Value2 = "dritte Unterliste"
$ws.Range("A27").NumberFormat = "@"
$ws.Range("B27").Value2 = "Spezialzeichen 1&2-%*_0 sind eingebettet"
$ws.Range("B27").NumberFormat = "@"
$ws.Range("D27").Value2 = "sehr"

# Row 28
$ws.Range("A28").Value2 = "dritte Unterliste"
$ws.Range("A28").NumberFormat = "@"
$ws.Range("B28").Value2 = "Spezialzeichen 1&2-%*_0 sind eingebettet"
$ws.Range("B28").NumberFormat = "@"
$ws.Range("D28").Value2 = "sehr"
$ws.Range("E28").Value2 = "tief"

# Row 29
$ws.Range("A29").Value2 = "dritte Unterliste"
$ws.Range("A29").NumberFormat = "@"
$ws.Range("B29").Value2 = "Spezialzeichen 1&2-%*_0 sind eingebettet"
$ws.Range("B29").NumberFormat = "@"
$ws.Range("D29").Value2 = "sehr"
$ws.Range("E29").Value2 = "tief"
$ws.Range("F29").Value2 = "verschachtelt!"

# Match row height used throughout the sheet's lower block.
$ws.Rows.Item(27).RowHeight = 16
$ws.Rows.Item(28).RowHeight = 16
$ws.Rows.Item(29).RowHeight = 16

# Move / resize the window and update the active selection like the
# recorded session did.
$ws.Range("D30").Select()
